$d = $word.ActiveDocument
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*tBbmSupihNo*") {
        $p.Range.Font.TextColor.ObjectThemeColor = 6
    }
}
